# Update stats for 2025-10 (row 23 of Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B23").Value = 6312
$ws.Range("C23").Value = 1001
$ws.Range("D23").Value = 5909622
$ws.Range("E23").Value = 936.2519011406844
$ws.Range("F23").Value = 8.304735758407688
$ws.Range("G23").Value = 4.162330905306977
$ws.Range("H23").Value = 26.62807298957577
